$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# Descripcion / Fecha_inicio / Fecha_fin / Estado columns one place to
# the right (B..E) and carries their widths/styles/content along with
# them automatically.
$ws.Columns("A").Insert()

# Populate the new column: a "Sede" selector with a sample value.
# A2 is written before A1 so the shared-string table picks up "Sede
# Chile" (index 5) ahead of "Nombre_sede" (index 6), matching the
# order in which the template author generated the file.
$ws.Range("A2").Value = "Sede Chile"
$ws.Range("A1").Value = "Nombre_sede"
$ws.Range("A1").Font.Bold = $true

# New column A width. The stored xlsx width (27.6640625 character
# units) isn't exactly reproducible through the ColumnWidth COM
# property (it quantizes to 1/6ths of a character); 26.8 is the COM
# input that lands on the closest achievable stored width
# (27.666666666666668, off by ~0.0026 from the source value).
$ws.Columns("A").ColumnWidth = 26.8

# Update the two sample dates (now in columns C/D after the shift).
$ws.Range("C2").Value = 45110
$ws.Range("D2").Value = 45117

# Move the active selection like the source workbook.
[void]$ws.Range("D6").Select()

Write-Host "edit applied"
